# Commit: "Change ppt: Add s in the title slide 1"
# The title text box on slide 1 ("Buddies Hub") gets an "s" appended,
# turning it into "Buddies Hubs".

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$title = $s.Shapes.Title
$title.TextFrame.TextRange.Text = "Buddies Hubs"
